$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 290 - this shifts rows 290..364 down to 291..365,
# preserving all existing data/formatting of those rows.
$ws.Rows.Item(290).Insert()

# Populate the newly inserted row 290 with the new record.
# Most fields mirror the surrounding rows (same market/category/etc.),
# only the date, volume, prices and $/Kg differ.
$ws.Cells.Item(290, 1).Value = 4
$ws.Cells.Item(290, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(290, 3).Value = "Los Lagos"
$ws.Cells.Item(290, 4).Value = 44932
$ws.Cells.Item(290, 5).Value = 10
$ws.Cells.Item(290, 6).Value = 100112043
$ws.Cells.Item(290, 7).Value = "Pepino ensalada"
$ws.Cells.Item(290, 8).Value = "Sin especificar"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 400
$ws.Cells.Item(290, 11).Value = 19000
$ws.Cells.Item(290, 12).Value = 20000
$ws.Cells.Item(290, 13).Value = 19500
$ws.Cells.Item(290, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(290, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(290, 16).Value = 325
$ws.Cells.Item(290, 17).Value = 60
$ws.Cells.Item(290, 18).Value = "Hortaliza"
